# Energy data revised and corrected.
# Populate the new "country data iea" (column C) values for the rows that
# were still missing their ISO-country-code CSV file name.
# The order below matches the order the values were (re)entered in the
# source workbook, which determines the order new strings are appended
# to the shared-strings table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value  = "MZ.csv"   # Mozambique
$ws.Range("C11").Value = "LY.csv"   # Libya
$ws.Range("C12").Value = "MX.csv"   # Mexico
$ws.Range("C14").Value = "TN.csv"   # Tunisia
$ws.Range("C20").Value = "PK.csv"   # Pakistan
$ws.Range("C5").Value  = "DZ.csv"   # Algeria
$ws.Range("C7").Value  = "CG.csv"   # Congo, Rep.
$ws.Range("C8").Value  = "GH.csv"   # Ghana
$ws.Range("C13").Value = "NG.csv"   # Nigeria
$ws.Range("C16").Value = "GB.csv"   # United Kingdom
$ws.Range("C17").Value = "US.csv"   # United States
$ws.Range("C18").Value = "VE.csv"   # Venezuela, RB
$ws.Range("C19").Value = "CI.csv"   # Cote d'Ivoire

# Match the final cursor/selection position left in the saved file.
$ws.Range("C21").Select()
